# Apply cryptos-list price/volume refresh (GitHub Actions data update).
# Numeric-looking Price (column D) values are written with a leading
# apostrophe so Excel keeps them as text (matching the sheet's existing
# inlineStr storage) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.417.17"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.652.50"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'597.01"
$ws.Range("D6").Value = "'157.96"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "2.651.11"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "'28.08"
$ws.Range("D15").Value = "3.136.66"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = "68.353.81"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "2.664.03"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'11.78"
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").Value = "'75.18"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'10.02"
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'573.97"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").Value = "'8.24"
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("D33").Value = "'1.41"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "'1.66"
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'160.72"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "'0.373"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").Value = "'5.34"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  -7.65%  "
$ws.Range("D46").Value = "'158.69"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'3.91"
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").Value = "'21.96"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +6.47%  "
$ws.Range("D51").Value = "'0.0784"
$ws.Range("E51").Value = "  -0.13%  "
